$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E11: change value from "Active" to "Activef"
$ws.Cells.Item(11, 5).Value = "Activef"

# E20: change from text "Active" to numeric 2
$ws.Cells.Item(20, 5).Value = 2

# New row 21
$ws.Cells.Item(21, 1).Value = 1015
$ws.Cells.Item(21, 2).Value = "test PRJ"
$ws.Cells.Item(21, 3).Value = 15
$ws.Cells.Item(21, 4).Value = "PRJ-564"
$ws.Cells.Item(21, 5).Value = "Active"

# Update selection to E19
$ws.Range("E19").Select()
